$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Item table header change -----------------------------------------
# Cell A1 holds rich text: "Defect " (bold, inherited from cell style) +
# "*" (bold, red) -> new wording "Defect Type " + "*" (bold, red).
$cell = $ws.Range("A1")

# Replace just the leading "Defect " run's text, keeping the trailing "*"
# run (and its formatting) intact.
$defectChars = $cell.Characters(1, 7)
$defectChars.Text = "Defect Type "

# Re-assert the bold/red formatting on the trailing "*" run, since
# rewriting the sibling run's text can reset it to the default font.
$starChar = $cell.Characters(13, 1)
$starChar.Font.Bold = $true
$starChar.Font.Color = 255

# --- Active selection moves from B6 to B1 -------------------------------
$ws.Range("B1").Select()
